$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 (copy formatting from an existing header cell, e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Per-row time_taken values (literal timestamps as captured originally)
$timestamps = @(
    "2021-10-05 10:52:11.632757",
    "2021-10-05 10:52:11.632769",
    "2021-10-05 10:52:11.632773",
    "2021-10-05 10:52:11.632776",
    "2021-10-05 10:52:11.632780",
    "2021-10-05 10:52:11.632783",
    "2021-10-05 10:52:11.632786",
    "2021-10-05 10:52:11.632789",
    "2021-10-05 10:52:11.632793",
    "2021-10-05 10:52:11.632796",
    "2021-10-05 10:52:11.632799",
    "2021-10-05 10:52:11.632802",
    "2021-10-05 10:52:11.632805",
    "2021-10-05 10:52:11.632808",
    "2021-10-05 10:52:11.632811",
    "2021-10-05 10:52:11.632814",
    "2021-10-05 10:52:11.632818",
    "2021-10-05 10:52:11.632821",
    "2021-10-05 10:52:11.632824",
    "2021-10-05 10:52:11.632827",
    "2021-10-05 10:52:11.632830",
    "2021-10-05 10:52:11.632833",
    "2021-10-05 10:52:11.632836",
    "2021-10-05 10:52:11.632839",
    "2021-10-05 10:52:11.632843"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
